$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 65. This shifts the existing rows 65-114
# (Champion Membrillo price records) down to rows 66-115, matching the
# weekly-update diff where a new observation is prepended to this block
# and every later record keeps its original data but moves down one row.
$ws.Rows(65).Insert()

# Populate the newly inserted row 65 with a fresh weekly record. It is a
# duplicate of the record that is now on row 66 (the same market basket:
# Feria Lagunitas de Puerto Montt / Los Lagos / Membrillo / Champion /
# Primera), except for a newer observation date (D65).
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44827
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100104
$ws.Range("H65").Value = "Frutos de pepita"
$ws.Range("I65").Value = 100104003
$ws.Range("J65").Value = "Membrillo"
$ws.Range("K65").Value = "Champion"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 14000
$ws.Range("O65").Value = 15000
$ws.Range("P65").Value = 14500
$ws.Range("Q65").Value = "$/caja 18 kilos granel"
$ws.Range("R65").Value = "Región de O'Higgins"
$ws.Range("S65").Value = 806
$ws.Range("T65").Value = 18
